$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 12 (Z05_B01_P01), shifting all
# subsequent rows down by one.
$ws.Rows.Item(12).Insert()

# Copy the formatting from the row now below the new blank row (row 13,
# which held the old row 12 content) so the new row matches the existing
# data-row style exactly.
$ws.Range("A13:D13").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)

# Populate the new row with the inserted postulate data.
$ws.Range("A12").Value = "Z04_B03_P01"
$ws.Range("B12").Value = "Z04_B03"
$ws.Range("C12").Value = "Mobilität in den Bildungsbiografien fördern"
$ws.Range("D12").Value = "XXXMobilität in den Bildungsbiografien fördern"
